$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.236.40'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '1.829.79'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.95'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6087'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.45%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07118'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2821'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.96'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07680'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.63%  '
$ws.Range('D12').Value = '1.816.14'
$ws.Range('E12').Value = '  -1.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.818'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6364'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001002'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('D16').Value = '2.067.90'
$ws.Range('E16').Value = '  -1.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.43'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.902'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.65%  '
$ws.Range('D19').Value = '29.198.06'
$ws.Range('E19').Value = '  -0.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '229.52'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.09%  '
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.026'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.000'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.39%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.093'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.76%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1290'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.67'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.62%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.511'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.77%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06553'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.456'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.851'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.844'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  -5.56%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6544'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.52%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.545'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.767'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.12%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '1.223.77'
$ws.Range('E39').Value = '  -1.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01761'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.633'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9267'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.11%  '
$ws.Range('D45').Value = '1.981.62'
$ws.Range('E45').Value = '  -1.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000117'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.613'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.83%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.569'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.533'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.24%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05547'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.32%  '
